{"js": "const replacements = [\n  [\"2025-07-10 Thursday\", \"2025-07-11 Friday\"],\n  [\"27\u00d716=432\", \"29\u00d712=348\"],\n  [\"57\u00d737=2109\", \"50\u00d771=3550\"],\n  [\"59\u00d793=5487\", \"89\u00d766=5874\"],\n  [\"88\u00d740=3520\", \"20\u00d719=380\"],\n  [\"76\u00d713=988\", \"77\u00d714=1078\"],\n  [\"47\u00d798=4606\", \"24\u00d756=1344\"],\n  [\"21\u00d734=714\", \"88\u00d762=5456\"],\n  [\"75\u00d792=6900\", \"60\u00d779=4740\"],\n  [\"11\u00d724=264\", \"71\u00d748=3408\"],\n  [\"23\u00d742=966\", \"55\u00d741=2255\"],\n  [\"13\u00d767=871\", \"69\u00d774=5106\"],\n  [\"44\u00d721=924\", \"25\u00d763=1575\"],\n  [\"90\u00d786=7740\", \"24\u00d718=432\"],\n  [\"48\u00d765=3120\", \"73\u00d790=6570\"],\n  [\"29\u00d754=1566\", \"61\u00d744=2684\"],\n  [\"97\u00d728=2716\", \"79\u00d743=3397\"],\n  [\"17\u00d754=918\", \"17\u00d770=1190\"],\n  [\"49\u00d762=3038\", \"63\u00d748=3024\"],\n  [\"68\u00d751=3468\", \"57\u00d732=1824\"],\n  [\"48\u00d737=1776\", \"57\u00d782=4674\"],\n  [\"88\u00d758=5104\", \"55\u00d715=825\"],\n  [\"70\u00d759=4130\", \"40\u00d777=3080\"],\n  [\"39\u00d766=2574\", \"48\u00d723=1104\"],\n  [\"32\u00d799=3168\", \"31\u00d753=1643\"],\n  [\"23\u00d752=1196\", \"28\u00d714=392\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\nreturn \"done\";", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-07-10 Thursday\", \"2025-07-11 Friday\"),\n    @(\"27\u00d716=432\", \"29\u00d712=348\"),\n    @(\"57\u00d737=2109\", \"50\u00d771=3550\"),\n    @(\"59\u00d793=5487\", \"89\u00d766=5874\"),\n    @(\"88\u00d740=3520\", \"20\u00d719=380\"),\n    @(\"76\u00d713=988\", \"77\u00d714=1078\"),\n    @(\"47\u00d798=4606\", \"24\u00d756=1344\"),\n    @(\"21\u00d734=714\", \"88\u00d762=5456\"),\n    @(\"75\u00d792=6900\", \"60\u00d779=4740\"),\n    @(\"11\u00d724=264\", \"71\u00d748=3408\"),\n    @(\"23\u00d742=966\", \"55\u00d741=2255\"),\n    @(\"13\u00d767=871\", \"69\u00d774=5106\"),\n    @(\"44\u00d721=924\", \"25\u00d763=1575\"),\n    @(\"90\u00d786=7740\", \"24\u00d718=432\"),\n    @(\"48\u00d765=3120\", \"73\u00d790=6570\"),\n    @(\"29\u00d754=1566\", \"61\u00d744=2684\"),\n    @(\"97\u00d728=2716\", \"79\u00d743=3397\"),\n    @(\"17\u00d754=918\", \"17\u00d770=1190\"),\n    @(\"49\u00d762=3038\", \"63\u00d748=3024\"),\n    @(\"68\u00d751=3468\", \"57\u00d732=1824\"),\n    @(\"48\u00d737=1776\", \"57\u00d782=4674\"),\n    @(\"88\u00d758=5104\", \"55\u00d715=825\"),\n    @(\"70\u00d759=4130\", \"40\u00d777=3080\"),\n    @(\"39\u00d766=2574\", \"48\u00d723=1104\"),\n    @(\"32\u00d799=3168\", \"31\u00d753=1643\"),\n    @(\"23\u00d752=1196\", \"28\u00d714=392\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output (\"NOT FOUND: \" + $oldText)\n    }\n}\n\nWrite-Output \"Done\""}
